$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3018.3333
$ws.Range("I106").Value = 2538.889
$ws.Range("J106").Value = 3497.7778
$ws.Range("K106").Value = 2538.889
$ws.Range("L106").Value = 3497.7778
$ws.Range("M106").Value = -1907.889
$ws.Range("N106").Value = -4759.7778

$ws.Range("H112").Value = 5615.5
$ws.Range("I112").Value = 742
$ws.Range("J112").Value = 7240
$ws.Range("K112").Value = 2226
$ws.Range("L112").Value = 21720
$ws.Range("M112").Value = -1118
$ws.Range("N112").Value = -23936

$ws.Range("H113").Value = 8141.811
$ws.Range("I113").Value = 16984.924
$ws.Range("J113").Value = 3351.7917
$ws.Range("K113").Value = 16984.924
$ws.Range("L113").Value = 3351.7917
$ws.Range("M113").Value = -13730.924
$ws.Range("N113").Value = -9859.7917

$ws.Range("H132").Value = 32076184
$ws.Range("I132").Value = 42453304
$ws.Range("K132").Value = 127359912
$ws.Range("M132").Value = -127357382

$ws.Range("H137").Value = 210088.61
$ws.Range("I137").Value = 314591.2
$ws.Range("J137").Value = 1083.4736
$ws.Range("K137").Value = 943773.6000000001
$ws.Range("L137").Value = 3250.4208
$ws.Range("M137").Value = -941223.6000000001
$ws.Range("N137").Value = -8350.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5280.905
$ws.Range("I61").Value = 5310.4736
$ws.Range("K61").Value = 5310.4736
$ws.Range("M61").Value = -5098.4736

$ws.Range("H74").Value = 1029.4166
$ws.Range("I74").Value = 759.3333
$ws.Range("J74").Value = 1299.5
$ws.Range("K74").Value = 759.3333
$ws.Range("L74").Value = 1299.5
$ws.Range("M74").Value = 114.6667
$ws.Range("N74").Value = -3047.5

$ws.Range("H77").Value = 1029.4166
$ws.Range("I77").Value = 759.3333
$ws.Range("J77").Value = 1299.5
$ws.Range("K77").Value = 3796.6665
$ws.Range("L77").Value = 6497.5
$ws.Range("M77").Value = 571.3334999999997
$ws.Range("N77").Value = -15233.5

$ws.Range("H132").Value = 3908129.8
$ws.Range("I132").Value = 6580768.5
$ws.Range("J132").Value = 1965.4615
$ws.Range("K132").Value = 19742305.5
$ws.Range("L132").Value = 5896.3845
$ws.Range("M132").Value = -19739775.5
$ws.Range("N132").Value = -10956.3845

$ws.Range("H136").Value = 5280.905
$ws.Range("I136").Value = 5310.4736
$ws.Range("K136").Value = 15931.4208
$ws.Range("M136").Value = -13381.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 83335304
$ws.Range("I107").Value = 100001976
$ws.Range("J107").Value = 1950
$ws.Range("K107").Value = 100001976
$ws.Range("L107").Value = 1950
$ws.Range("M107").Value = -100000056
$ws.Range("N107").Value = -5790

$ws.Range("H134").Value = 9023557
$ws.Range("I134").Value = 10769661
$ws.Range("J134").Value = 2019
$ws.Range("K134").Value = 32308983
$ws.Range("L134").Value = 6057
$ws.Range("M134").Value = -32306448
$ws.Range("N134").Value = -11127

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6644.327
$ws.Range("I31").Value = 1063.5454
$ws.Range("J31").Value = 28967.455
$ws.Range("K31").Value = 1063.5454
$ws.Range("L31").Value = 28967.455
$ws.Range("M31").Value = -768.5454
$ws.Range("N31").Value = -29557.455

$ws.Range("H34").Value = 6644.327
$ws.Range("I34").Value = 1063.5454
$ws.Range("J34").Value = 28967.455
$ws.Range("K34").Value = 1063.5454
$ws.Range("L34").Value = 28967.455
$ws.Range("M34").Value = -861.5454
$ws.Range("N34").Value = -29371.455

$ws.Range("H58").Value = 2526513.2
$ws.Range("I58").Value = 3270272.5
$ws.Range("J58").Value = 9174.154
$ws.Range("K58").Value = 3270272.5
$ws.Range("L58").Value = 9174.154
$ws.Range("M58").Value = -3270069.5
$ws.Range("N58").Value = -9580.154

$ws.Range("H99").Value = 76925330
$ws.Range("I99").Value = 76925330
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 76925330
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -76923832
$ws.Range("N99").ClearContents()

$ws.Range("H105").Value = 1449.1875
$ws.Range("I105").Value = 1414.3846
$ws.Range("J105").Value = 1600
$ws.Range("K105").Value = 1414.3846
$ws.Range("L105").Value = 1600
$ws.Range("M105").Value = 332.6153999999999
$ws.Range("N105").Value = -5094

$ws.Range("H126").Value = 76925330
$ws.Range("I126").Value = 76925330
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 230775990
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -230773520
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 6806869
$ws.Range("I132").Value = 10754651
$ws.Range("J132").Value = 7911.1113
$ws.Range("K132").Value = 32263953
$ws.Range("L132").Value = 23733.3339
$ws.Range("M132").Value = -32261423
$ws.Range("N132").Value = -28793.3339

$ws.Range("H134").Value = 29167990
$ws.Range("I134").Value = 40323720
$ws.Range("K134").Value = 120971160
$ws.Range("M134").Value = -120968625

$ws.Range("H136").Value = 2526513.2
$ws.Range("I136").Value = 3270272.5
$ws.Range("J136").Value = 9174.154
$ws.Range("K136").Value = 9810817.5
$ws.Range("L136").Value = 27522.462
$ws.Range("M136").Value = -9808267.5
$ws.Range("N136").Value = -32622.462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2587.3333
$ws.Range("I123").Value = 1225.4286
$ws.Range("J123").Value = 4494
$ws.Range("K123").Value = 3676.2858
$ws.Range("L123").Value = 13482
$ws.Range("M123").Value = -1226.2858
$ws.Range("N123").Value = -18382

$ws.Range("H129").Value = 1831.2858
$ws.Range("I129").Value = 680
$ws.Range("J129").Value = 2145.2727
$ws.Range("K129").Value = 2040
$ws.Range("L129").Value = 6435.8181
$ws.Range("M129").Value = 2960
$ws.Range("N129").Value = -16435.8181

$ws.Range("H130").Value = 12213.333
$ws.Range("I130").Value = 25605
$ws.Range("K130").Value = 76815
$ws.Range("M130").Value = -71795

$ws.Range("J131").Value = 1018322.4
$ws.Range("L131").Value = 3054967.2
$ws.Range("N131").Value = -3065047.2

$ws.Range("H133").Value = 7987.9644
$ws.Range("I133").Value = 1714.1428
$ws.Range("J133").Value = 8558.312
$ws.Range("K133").Value = 5142.428400000001
$ws.Range("L133").Value = 25674.936
$ws.Range("M133").Value = -82.42840000000069
$ws.Range("N133").Value = -35794.936

$ws.Range("H134").Value = 4845.143
$ws.Range("I134").Value = 4837.4
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 14512.2
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -9442.199999999999
$ws.Range("N134").Value = -25140

$ws.Range("H136").Value = 2690
$ws.Range("I136").Value = 1207.5
$ws.Range("J136").Value = 4666.6665
$ws.Range("K136").Value = 3622.5
$ws.Range("L136").Value = 13999.9995
$ws.Range("M136").Value = 1477.5
$ws.Range("N136").Value = -24199.9995

$ws.Range("H137").Value = 23859948
$ws.Range("I137").Value = 74199.92999999999
$ws.Range("J137").Value = 71431440
$ws.Range("K137").Value = 222599.79
$ws.Range("L137").Value = 214294320
$ws.Range("M137").Value = -217499.79
$ws.Range("N137").Value = -214304520

$ws.Range("H138").Value = 3449.3125
$ws.Range("I138").Value = 3724.0833
$ws.Range("J138").Value = 2625
$ws.Range("K138").Value = 11172.2499
$ws.Range("L138").Value = 7875
$ws.Range("M138").Value = -6032.249899999999
$ws.Range("N138").Value = -18155

$ws.Range("H139").Value = 984.61536
$ws.Range("I139").Value = 816.6667
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 2450.0001
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 2689.9999
$ws.Range("N139").Value = -19280

$ws.Range("H140").Value = 1147.8125
$ws.Range("I140").Value = 1091
$ws.Range("K140").Value = 3273
$ws.Range("M140").Value = 1907

$ws.Range("H141").Value = 2665.7917
$ws.Range("I141").Value = 1619.4762
$ws.Range("K141").Value = 4858.4286
$ws.Range("M141").Value = 321.5713999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 50003820
$ws.Range("I132").Value = 86958296
$ws.Range("J132").Value = 6580.5293
$ws.Range("K132").Value = 260874888
$ws.Range("L132").Value = 19741.5879
$ws.Range("M132").Value = -260872358
$ws.Range("N132").Value = -24801.5879

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4349751
$ws.Range("I132").Value = 6452798.5
$ws.Range("J132").Value = 3451.8667
$ws.Range("K132").Value = 19358395.5
$ws.Range("L132").Value = 10355.6001
$ws.Range("M132").Value = -19355865.5
$ws.Range("N132").Value = -15415.6001

$ws.Range("H136").Value = 3072.0615
$ws.Range("I136").Value = 3240.3518
$ws.Range("J136").Value = 2245.9092
$ws.Range("K136").Value = 9721.055399999999
$ws.Range("L136").Value = 6737.7276
$ws.Range("M136").Value = -7171.055399999999
$ws.Range("N136").Value = -11837.7276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7762232.5
$ws.Range("I132").Value = 4762972
$ws.Range("J132").Value = 10281611
$ws.Range("K132").Value = 14288916
$ws.Range("L132").Value = 30844833
$ws.Range("M132").Value = -14286386
$ws.Range("N132").Value = -30849893

$ws.Range("H136").Value = 13186782
$ws.Range("I136").Value = 7288150
$ws.Range("J136").Value = 31251340
$ws.Range("K136").Value = 21864450
$ws.Range("L136").Value = 93754020
$ws.Range("M136").Value = -21861900
$ws.Range("N136").Value = -93759120
